$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2026-02-12 Thursday" "2026-02-13 Friday"

Replace-Text "13÷7=" "68÷4="
Replace-Text "36÷9=" "64÷5="
Replace-Text "34÷5=" "49÷6="
Replace-Text "21÷2=" "88÷5="
Replace-Text "36÷4=" "30÷7="

Replace-Text "42÷9=" "32÷4="
Replace-Text "65÷5=" "69÷4="
Replace-Text "58÷4=" "40÷6="
Replace-Text "62÷6=" "38÷8="
Replace-Text "17÷4=" "64÷3="

Replace-Text "46÷9=" "98÷8="
Replace-Text "42÷6=" "75÷4="
Replace-Text "72÷4=" "77÷7="
Replace-Text "72÷8=" "27÷6="
Replace-Text "62÷9=" "93÷2="

Replace-Text "85÷5=" "51÷5="
Replace-Text "17÷8=" "92÷6="
Replace-Text "98÷2=" "93÷6="
Replace-Text "81÷9=" "43÷7="
Replace-Text "24÷3=" "85÷8="

Replace-Text "58÷6=" "35÷8="
Replace-Text "97÷9=" "88÷8="
Replace-Text "92÷5=" "89÷2="
Replace-Text "14÷3=" "74÷4="
Replace-Text "78÷8=" "67÷5="
